$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Vacant/Occupied" occupancy-tracking section below the
# existing DoorState table (rows 5-7), matching the order the strings
# were originally typed in (col B down, then C7 before C5/C6, then D down).
$ws.Range("A5").Value = "Vacant?"
$ws.Range("B5").Value = "Vacant"
$ws.Range("B6").Value = "Occupied"
$ws.Range("B7").Value = "Occupants"
$ws.Range("C7").Value = "Everytime sensor 1 then sensor 2 add 1 and subtract for opposite"
$ws.Range("C5").Value = "if occupants is > 0 it is true"
$ws.Range("C6").Value = "If Occupants is < 0 it is true"
$ws.Range("D5").Value = "Recognizes if a room is occupied"
$ws.Range("D6").Value = "Recognizes if a room is vacant"
$ws.Range("D7").Value = "Keeps track of the people in the room"


# Column C now holds the long "Logic" descriptions, so widen it to fit.
$ws.Columns.Item(3).ColumnWidth = 58.67

# Cursor ends up on the next empty row after the data entry.
$ws.Range("D8").Select()
